$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 105265700
$ws.Range("I62").Value = 71430130
$ws.Range("J62").Value = 200005280
$ws.Range("K62").Value = 71430130
$ws.Range("L62").Value = 200005280
$ws.Range("M62").Value = -71429506
$ws.Range("N62").Value = -200006528
$ws.Range("H64").Value = 3167.4583
$ws.Range("I64").Value = 3272.8572
$ws.Range("J64").Value = 3124.0588
$ws.Range("K64").Value = 3272.8572
$ws.Range("L64").Value = 3124.0588
$ws.Range("M64").Value = -3024.8572
$ws.Range("N64").Value = -3620.0588
$ws.Range("H65").Value = 105265700
$ws.Range("I65").Value = 71430130
$ws.Range("J65").Value = 200005280
$ws.Range("K65").Value = 357150650
$ws.Range("L65").Value = 1000026400
$ws.Range("M65").Value = -357147530
$ws.Range("N65").Value = -1000032640
$ws.Range("H67").Value = 3167.4583
$ws.Range("I67").Value = 3272.8572
$ws.Range("J67").Value = 3124.0588
$ws.Range("K67").Value = 3272.8572
$ws.Range("L67").Value = 3124.0588
$ws.Range("M67").Value = -2414.8572
$ws.Range("N67").Value = -4840.0588

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 142860030
$ws.Range("J88").Value = 200003310
$ws.Range("L88").Value = 200003310
$ws.Range("N88").Value = -200004122
$ws.Range("H91").Value = 142860030
$ws.Range("J91").Value = 200003310
$ws.Range("L91").Value = 200003310
$ws.Range("N91").Value = -200006118
$ws.Range("H97").Value = 1345.3226
$ws.Range("I97").Value = 1734.5217
$ws.Range("J97").Value = 226.375
$ws.Range("K97").Value = 1734.5217
$ws.Range("L97").Value = 226.375
$ws.Range("M97").Value = -1238.5217
$ws.Range("N97").Value = -1218.375
$ws.Range("H102").Value = 2547.1428
$ws.Range("I102").Value = 2388.3333
$ws.Range("J102").Value = 3500
$ws.Range("K102").Value = 2388.3333
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = -766.3332999999998
$ws.Range("N102").Value = -6744

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 115.75
$ws.Range("I80").Value = 50
$ws.Range("J80").Value = 137.66667
$ws.Range("K80").Value = 50
$ws.Range("L80").Value = 137.66667
$ws.Range("M80").Value = 948
$ws.Range("N80").Value = -2133.66667
$ws.Range("H83").Value = 115.75
$ws.Range("I83").Value = 50
$ws.Range("J83").Value = 137.66667
$ws.Range("K83").Value = 250
$ws.Range("L83").Value = 688.3333500000001
$ws.Range("M83").Value = 4742
$ws.Range("N83").Value = -10672.33335
$ws.Range("H86").Value = 11113084
$ws.Range("J86").Value = 1532.8572
$ws.Range("L86").Value = 1532.8572
$ws.Range("N86").Value = -3778.8572
$ws.Range("H89").Value = 11113084
$ws.Range("J89").Value = 1532.8572
$ws.Range("L89").Value = 7664.286
$ws.Range("N89").Value = -18896.286
$ws.Range("H94").Value = 8567.879999999999
$ws.Range("I94").Value = 391.8421
$ws.Range("J94").Value = 34458.668
$ws.Range("K94").Value = 391.8421
$ws.Range("L94").Value = 34458.668
$ws.Range("M94").Value = 59.15789999999998
$ws.Range("N94").Value = -35360.668
$ws.Range("H99").Value = 1139.4445
$ws.Range("I99").Value = 1147
$ws.Range("J99").Value = 1011
$ws.Range("K99").Value = 1147
$ws.Range("L99").Value = 1011
$ws.Range("M99").Value = 351
$ws.Range("N99").Value = -4007
$ws.Range("H105").Value = 2589.65
$ws.Range("I105").Value = 2167.3572
$ws.Range("J105").Value = 3575
$ws.Range("K105").Value = 2167.3572
$ws.Range("L105").Value = 3575
$ws.Range("M105").Value = -420.3571999999999
$ws.Range("N105").Value = -7069

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8719560
$ws.Range("I31").Value = 10579815
$ws.Range("J31").Value = 7410491.5
$ws.Range("K31").Value = 10579815
$ws.Range("L31").Value = 7410491.5
$ws.Range("M31").Value = -10579520
$ws.Range("N31").Value = -7411081.5
$ws.Range("H34").Value = 8719560
$ws.Range("I34").Value = 10579815
$ws.Range("J34").Value = 7410491.5
$ws.Range("K34").Value = 10579815
$ws.Range("L34").Value = 7410491.5
$ws.Range("M34").Value = -10579613
$ws.Range("N34").Value = -7410895.5
$ws.Range("H86").Value = 627490.1
$ws.Range("I86").Value = 1002122.6
$ws.Range("J86").Value = 3102.6667
$ws.Range("K86").Value = 1002122.6
$ws.Range("L86").Value = 3102.6667
$ws.Range("M86").Value = -1000999.6
$ws.Range("N86").Value = -5348.6667
$ws.Range("H89").Value = 627490.1
$ws.Range("I89").Value = 1002122.6
$ws.Range("J89").Value = 3102.6667
$ws.Range("K89").Value = 5010613
$ws.Range("L89").Value = 15513.3335
$ws.Range("M89").Value = -5004997
$ws.Range("N89").Value = -26745.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 94.5
$ws.Range("I23").Value = 80
$ws.Range("J23").Value = 96.57143000000001
$ws.Range("K23").Value = 240
$ws.Range("L23").Value = 289.71429
$ws.Range("M23").Value = -5
$ws.Range("N23").Value = -759.71429
$ws.Range("H97").Value = 5103132
$ws.Range("I97").Value = 11904990
$ws.Range("J97").Value = 1738.5
$ws.Range("K97").Value = 35714970
$ws.Range("L97").Value = 5215.5
$ws.Range("M97").Value = -35714474
$ws.Range("N97").Value = -6207.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 25003752
$ws.Range("I20").Value = 50000004
$ws.Range("K20").Value = 50000004
$ws.Range("M20").Value = -49999759
$ws.Range("H80").Value = 2835.25
$ws.Range("I80").Value = 2972.8572
$ws.Range("J80").Value = 2806.0605
$ws.Range("K80").Value = 2972.8572
$ws.Range("L80").Value = 2806.0605
$ws.Range("M80").Value = -1974.8572
$ws.Range("N80").Value = -4802.0605
$ws.Range("H83").Value = 2835.25
$ws.Range("I83").Value = 2972.8572
$ws.Range("J83").Value = 2806.0605
$ws.Range("K83").Value = 14864.286
$ws.Range("L83").Value = 14030.3025
$ws.Range("M83").Value = -9872.286
$ws.Range("N83").Value = -24014.3025
$ws.Range("H97").Value = 571.3158
$ws.Range("I97").Value = 722.0833
$ws.Range("J97").Value = 312.85715
$ws.Range("K97").Value = 722.0833
$ws.Range("L97").Value = 312.85715
$ws.Range("M97").Value = -226.0833
$ws.Range("N97").Value = -1304.85715
$ws.Range("H113").Value = 4078
$ws.Range("I113").Value = 1695
$ws.Range("J113").Value = 5666.6665
$ws.Range("K113").Value = 1695
$ws.Range("L113").Value = 5666.6665
$ws.Range("M113").Value = 475
$ws.Range("N113").Value = -10006.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6797.4287
$ws.Range("I62").Value = 4701
$ws.Range("J62").Value = 7636
$ws.Range("K62").Value = 4701
$ws.Range("L62").Value = 7636
$ws.Range("M62").Value = -4077
$ws.Range("N62").Value = -8884
$ws.Range("H65").Value = 6797.4287
$ws.Range("I65").Value = 4701
$ws.Range("J65").Value = 7636
$ws.Range("K65").Value = 23505
$ws.Range("L65").Value = 38180
$ws.Range("M65").Value = -20385
$ws.Range("N65").Value = -44420

Write-Host "Updated 181 cells across ALC, ARM, BSM, CRP, CUL, GSM, WVR sheets"
